# Reorder the worksheets: "review_info" moves before "hotel_info" so the
# tab order becomes review_info, hotel_info (names/tab positions swap while
# keeping each sheet's own data attached to its name).
$wb = $excel.ActiveWorkbook
$wsHotelBeforeMove = $wb.Worksheets.Item("hotel_info")
$wsReview = $wb.Worksheets.Item("review_info")
$wsReview.Move($wsHotelBeforeMove)

# Re-fetch the "hotel_info" worksheet by name now that the move has
# happened, since worksheet references resolve by position, not identity.
$wsHotel = $wb.Worksheets.Item("hotel_info")

# Insert a new "State" column into hotel_info, between Hotel_Name and City.
$wsHotel.Columns.Item(3).Insert()
$wsHotel.Cells.Item(1, 3).Value = "State"
$wsHotel.Cells.Item(2, 3).Value = "Louisiana"
